$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell 'D2' '34.953.37'
Set-TextCell 'E2' '  +0.67%  '

Set-TextCell 'D3' '1.840.36'
Set-TextCell 'E3' '  +1.84%  '

Set-TextCell 'E4' '  +0.07%  '

Set-TextCell 'D5' '232.41'
Set-TextCell 'E5' '  +0.52%  '

Set-TextCell 'D6' '0.618'
Set-TextCell 'E6' '  +2.21%  '

Set-TextCell 'E7' '  +0.04%  '

Set-TextCell 'D8' '41.12'
Set-TextCell 'E8' '  +4.75%  '

Set-TextCell 'E9' '  +2.89%  '

Set-TextCell 'D10' '0.0692'
Set-TextCell 'E10' '  +1.75%  '

Set-TextCell 'E11' '  -1.05%  '

Set-TextCell 'D12' '2.106.72'
Set-TextCell 'E12' '  +1.89%  '

Set-TextCell 'D13' '11.36'
Set-TextCell 'E13' '  +4.01%  '

Set-TextCell 'D14' '1.831.11'
Set-TextCell 'E14' '  +1.73%  '

Set-TextCell 'D15' '0.670'
Set-TextCell 'E15' '  +1.27%  '

Set-TextCell 'E16' '  +2.20%  '

Set-TextCell 'D17' '34.934.19'
Set-TextCell 'E17' '  +0.64%  '

Set-TextCell 'D18' '69.84'
Set-TextCell 'E18' '  +0.34%  '

Set-TextCell 'D19' '0.0₃0790'
Set-TextCell 'E19' '  +0.86%  '

Set-TextCell 'D20' '239.83'
Set-TextCell 'E20' '  +0.09%  '

Set-TextCell 'E21' '  +2.73%  '

Set-TextCell 'D22' '4.76'
Set-TextCell 'E22' '  +1.88%  '

Set-TextCell 'E23' '  +0.04%  '

Set-TextCell 'E24' '  +1.00%  '

Set-TextCell 'D25' '172.02'
Set-TextCell 'E25' '  -0.31%  '

Set-TextCell 'E26' '  +1.93%  '

Set-TextCell 'D27' '17.44'
Set-TextCell 'E27' '  +1.56%  '

Set-TextCell 'E28' '  +3.64%  '

Set-TextCell 'D29' '1.70'
Set-TextCell 'E29' '  +11.53%  '

Set-TextCell 'E30' '  +0.09%  '

Set-TextCell 'E31' '  +0.92%  '

Set-TextCell 'E32' '  -1.09%  '

Set-TextCell 'E33' '  -0.99%  '

Set-TextCell 'E34' '  +22.78%  '

Set-TextCell 'D35' '1.95'
Set-TextCell 'E35' '  +10.48%  '

Set-TextCell 'D36' '1.26'
Set-TextCell 'E36' '  +1.91%  '

Set-TextCell 'D37' '0.756'
Set-TextCell 'E37' '  +8.12%  '

Set-TextCell 'E38' '  +10.19%  '

Set-TextCell 'D39' '89.76'
Set-TextCell 'E39' '  -1.74%  '

Set-TextCell 'E40' '  +3.28%  '

Set-TextCell 'D41' '1.336.74'
Set-TextCell 'E41' '  +1.94%  '

Set-TextCell 'D42' '14.58'
Set-TextCell 'E42' '  +2.59%  '

Set-TextCell 'D43' '2.40'
Set-TextCell 'E43' '  -2.18%  '

Set-TextCell 'E44' '  +1.49%  '

Set-TextCell 'E45' '  +3.43%  '

Set-TextCell 'D46' '0.0529'
Set-TextCell 'E46' '  +3.55%  '

Set-TextCell 'D47' '6.32'
Set-TextCell 'E47' '  +2.87%  '

Set-TextCell 'D48' '2.022.90'
Set-TextCell 'E48' '  +1.39%  '

Set-TextCell 'D49' '10.97'
Set-TextCell 'E49' '  +64.14%  '

Set-TextCell 'E50' '  +0.09%  '

Set-TextCell 'B51' 'Cronos'
Set-TextCell 'C51' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextCell 'D51' '0.0668'
Set-TextCell 'E51' '  -0.54%  '
